$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added; insert a row for it above the existing row 18
# (pushing all subsequent records down by one) and populate it with the new data.
$ws.Rows.Item(18).Insert()

$ws.Cells.Item(18, 1).Value = 10
$ws.Cells.Item(18, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(18, 3).Value = "La Araucanía"
$ws.Cells.Item(18, 4).Value = 44466
$ws.Cells.Item(18, 5).Value = 9
$ws.Cells.Item(18, 6).Value = 100112026
$ws.Cells.Item(18, 7).Value = "Haba"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 80
$ws.Cells.Item(18, 11).Value = 14000
$ws.Cells.Item(18, 12).Value = 14000
$ws.Cells.Item(18, 13).Value = 14000
$ws.Cells.Item(18, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(18, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(18, 16).Value = 560
$ws.Cells.Item(18, 17).Value = 25
$ws.Cells.Item(18, 18).Value = "Hortaliza"
